# Fruta / hortaliza, semanal
# Insert a new weekly record at row 805 in the Kiwi price sheet, shifting
# all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(805).Insert()

$ws.Range("A805").Value = 10
$ws.Range("B805").Value = "Vega Modelo de Temuco"
$ws.Range("C805").Value = "La Araucanía"
$ws.Range("D805").Value = 45223
$ws.Range("E805").Value = 9
$ws.Range("F805").Value = "Fruta"
$ws.Range("G805").Value = 100101
$ws.Range("H805").Value = "Berries"
$ws.Range("I805").Value = 100101007
$ws.Range("J805").Value = "Kiwi"
$ws.Range("K805").Value = "Hayward"
$ws.Range("L805").Value = "Primera"
$ws.Range("M805").Value = 210
$ws.Range("N805").Value = 35000
$ws.Range("O805").Value = 35000
$ws.Range("P805").Value = 35000
$ws.Range("Q805").Value = "`$/caja 15 kilos"
$ws.Range("R805").Value = "Región de O'Higgins"
$ws.Range("S805").Value = 2333
$ws.Range("T805").Value = 15
